$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List_ID")

$ws.Range("C9").Value = "23/07/1992"
$ws.Range("D9").Value = "xkknnasnpqpxdzj@gmail.com"
$ws.Range("E9").Value = "pmpbuTRSZW5"
$ws.Range("F9").Value = "pass"

$ws.Range("C10").Value = "19/12/1990"
$ws.Range("D10").Value = "epmxnqagoyunrzs@gmail.com"
$ws.Range("E10").Value = "bjnfdMIOHT5"
$ws.Range("F10").Value = "pass"
